$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table refresh (GitHub Actions scheduled update).
# Column layout: A=rank(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h).
# Some D-column values look numeric (e.g. "0.618") but must be stored as literal
# text, matching every other cell in the column; pre-setting NumberFormat to "@"
# (Text) before the assignment stops Excel from auto-coercing them to numbers.

# Row 2
$ws.Cells.Item(2, 4).Value = "43.606.12"
$ws.Cells.Item(2, 5).Value = "  +1.86%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.330.98"
$ws.Cells.Item(3, 5).Value = "  +1.68%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

# Row 5
$ws.Cells.Item(5, 2).Value = "Solana"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "109.29"
$ws.Cells.Item(5, 5).Value = "  +5.13%  "

# Row 6
$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "311.33"
$ws.Cells.Item(6, 5).Value = "  -1.48%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +1.13%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.11%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.618"
$ws.Cells.Item(9, 5).Value = "  +2.88%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +4.33%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0920"
$ws.Cells.Item(11, 5).Value = "  +1.85%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +1.89%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.74%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -0.24%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.52"
$ws.Cells.Item(15, 5).Value = "  +1.55%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.683.14"
$ws.Cells.Item(16, 5).Value = "  +1.60%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.330.82"
$ws.Cells.Item(17, 5).Value = "  +0.88%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "43.489.66"
$ws.Cells.Item(18, 5).Value = "  +1.77%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.58"
$ws.Cells.Item(19, 5).Value = "  +1.84%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +1.64%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.12"
$ws.Cells.Item(21, 5).Value = "  -4.98%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "74.26"
$ws.Cells.Item(22, 5).Value = "  +0.39%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -1.36%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "269.37"
$ws.Cells.Item(24, 5).Value = "  +2.34%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +3.23%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.00"
$ws.Cells.Item(26, 5).Value = "  -0.17%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.59"
$ws.Cells.Item(27, 5).Value = "  +7.88%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "11.19"
$ws.Cells.Item(28, 5).Value = "  +2.93%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.29"
$ws.Cells.Item(29, 5).Value = "  -2.76%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "38.94"
$ws.Cells.Item(30, 5).Value = "  +4.14%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "22.66"
$ws.Cells.Item(31, 5).Value = "  +1.39%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "167.35"
$ws.Cells.Item(32, 5).Value = "  +0.51%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0888"
$ws.Cells.Item(33, 5).Value = "  +2.06%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.84"
$ws.Cells.Item(34, 5).Value = "  +9.96%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.26%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.77"
$ws.Cells.Item(36, 5).Value = "  +4.55%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -1.76%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0365"
$ws.Cells.Item(38, 5).Value = "  +4.44%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.83"
$ws.Cells.Item(39, 5).Value = "  +0.50%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +6.15%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.71"
$ws.Cells.Item(41, 5).Value = "  +9.16%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "105.20"
$ws.Cells.Item(42, 5).Value = "  +14.10%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "71.79"
$ws.Cells.Item(43, 5).Value = "  +3.52%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Celestia"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "13.41"
$ws.Cells.Item(44, 5).Value = "  +10.64%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.236"
$ws.Cells.Item(45, 5).Value = "  +2.86%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.06%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "114.28"
$ws.Cells.Item(47, 5).Value = "  +0.32%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "1.672.71"
$ws.Cells.Item(48, 5).Value = "  -2.87%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "ordi"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "76.78"
$ws.Cells.Item(49, 5).Value = "  -4.27%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.36"
$ws.Cells.Item(50, 5).Value = "  +4.21%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "FraxShare"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "8.95"
$ws.Cells.Item(51, 5).Value = "  +2.19%  "
